$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Validation_accuracy) and D (Training_accuracy) new values for rows 2-26
$values = @{
    2  = @(0.70703125, 0.6849056482315063)
    3  = @(0.88671875, 0.9448113441467285)
    4  = @(0.91015625, 0.9886792302131653)
    5  = @(0.91796875, 0.9990565776824951)
    6  = @(0.91796875, 1)
    7  = @(0.91796875, 1)
    8  = @(0.91796875, 1)
    9  = @(0.91796875, 1)
    10 = @(0.91796875, 1)
    11 = @(0.91796875, 1)
    12 = @(0.91796875, 1)
    13 = @(0.91796875, 1)
    14 = @(0.91796875, 1)
    15 = @(0.91796875, 1)
    16 = @(0.91796875, 1)
    17 = @(0.91796875, 1)
    18 = @(0.91796875, 1)
    19 = @(0.91796875, 1)
    20 = @(0.91796875, 1)
    21 = @(0.91796875, 1)
    22 = @(0.9140625, 1)
    23 = @(0.9140625, 1)
    24 = @(0.9140625, 1)
    25 = @(0.9140625, 1)
    26 = @(0.921875, 1)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 3).Value = $pair[0]
    $ws.Cells.Item($row, 4).Value = $pair[1]
}
